$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.556.27"
$ws.Range("E2").Value = "  +1.81%  "
$ws.Range("D3").Value = "3.567.51"
$ws.Range("E3").Value = "  +1.99%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "620.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.92%  "
$ws.Range("D7").Value = "3.567.52"
$ws.Range("E7").Value = "  +2.02%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +2.27%  "
$ws.Range("E10").Value = "  +5.67%  "
$ws.Range("E11").Value = "  +5.68%  "
$ws.Range("E12").Value = "  +4.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000221"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.34%  "
$ws.Range("D15").Value = "4.175.71"
$ws.Range("E15").Value = "  +2.18%  "
$ws.Range("D16").Value = "3.564.38"
$ws.Range("E16").Value = "  +2.04%  "
$ws.Range("D17").Value = "68.475.08"
$ws.Range("E17").Value = "  +1.82%  "
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.07%  "
$ws.Range("E20").Value = "  +6.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +10.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "455.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.643"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.11%  "
$ws.Range("E25").Value = "  +1.49%  "
$ws.Range("D26").Value = "3.712.19"
$ws.Range("E26").Value = "  +2.14%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("E28").Value = "  +4.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.15"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +10.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.19%  "
$ws.Range("E31").Value = "  +3.53%  "
$ws.Range("E32").Value = "  +5.20%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.36"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.56%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "26.18"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.93%  "
$ws.Range("E36").Value = "  +3.56%  "
$ws.Range("D37").Value = "3.565.19"
$ws.Range("E37").Value = "  +2.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.27"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.40%  "
$ws.Range("E39").Value = "  +9.49%  "
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "179.15"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.39%  "
$ws.Range("E42").Value = "  +5.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.60"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "31.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +15.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.898"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.93%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.71%  "
$ws.Range("E48").Value = "  +5.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.67"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.263"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.73%  "
